# "cambios a la tabla" - add columns F and G to the table, and adjust
# the border of column E so the box stretches into the new columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Fill the new F and G columns (rows 2-10) with the same boxed
#        border style already used by columns C/D (cellXfs style index 1),
#        obtained here by copying the format from an existing boxed cell.
$ws.Range("C2").Copy()
$ws.Range("F2:G10").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 6).Value = 222   # column F
    $ws.Cells.Item($r, 7).Value = 225   # column G
}

# --- 2. Column E (E2:E9) loses its right border, since the boxed area now
#        continues on into the new F/G columns instead of ending at E.
$eRange = $ws.Range("E2:E9")
$eRange.Borders.Item(10).LineStyle = -4142  # xlEdgeRight -> none

# --- 3. Update the selection shown in the sheet view to match the new
#        active cell recorded after the edit.
$ws.Range("I8").Select()
